# Applies the two text edits described by the diff:
#  1. Slide 3 ("Что нужно, чтобы использовать сервис"), bullet shape "Объект 2",
#     2nd paragraph: "Телефон с " / "доступом к Wi" / "-Fi"
#                 -> "Телефон с доступом к " / "Wi" / "-Fi"
#  2. Slide 6 ("Сам сайт"), title shape "Заголовок 1":
#     "Сам сайт" -> "Наше устройство"

$p = $ppt.ActivePresentation

# --- Edit 1: slide 3, "Объект 2" bullet list, paragraph 2 ---------------
$slide3 = $p.Slides.Item(3)
$bulletShape = $slide3.Shapes.Item("Объект 2")
$para2 = $bulletShape.TextFrame.TextRange.Paragraphs(2, 1)

$run1 = $para2.Runs(1, 1)
$run1.Text = "Телефон с доступом к "

$run2 = $para2.Runs(2, 1)
$run2.Text = "Wi"

# 3rd run ("-Fi") is unchanged.

# --- Edit 2: slide 6, title shape "Заголовок 1" --------------------------
$slide6 = $p.Slides.Item(6)
$titleShape = $slide6.Shapes.Item("Заголовок 1")
$titleShape.TextFrame.TextRange.Text = "Наше устройство"
